$d = $word.ActiveDocument
$cr = [char]13

# 1. Update the letter date: "September 19, 2025" -> "September 21, 2025"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $text = $p.Range.Text.TrimEnd($cr)
    if ($text -eq "September 19, 2025") {
        $p.Range.Text = "September 21, 2025"
        break
    }
}

# 2. Split the mailing-address paragraph "2933 Lamory Pl, Santa Clara CA 95051"
#    into two paragraphs: "2933 Lamory Pl" and "Santa Clara, CA 95051"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $text = $p.Range.Text.TrimEnd($cr)
    if ($text -eq "2933 Lamory Pl, Santa Clara CA 95051") {
        $p.Range.Text = "2933 Lamory Pl${cr}Santa Clara, CA 95051"
        break
    }
}

# 3. Remove the empty "No Spacing" paragraph that immediately follows the
#    "Board of Directors" signature line.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $text = $p.Range.Text.TrimEnd($cr)
    if ($text -like "*Board of Directors") {
        $d.Paragraphs($i + 1).Range.Delete()
        break
    }
}
